$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values ---

# Row 5: add new E5 value
$ws.Range("E5").Value = "FLORES"

# Row 6: CP005_Empresa
$ws.Range("A6").Value = "CP005_Empresa"
$ws.Range("B6").Value = "Centro Banca Empresas`nAsesoramiento especializado en todo el país, ese es nuestro centro."

# Row 7: CP006_Cajero
$ws.Range("A7").Value = "CP006_Cajero"
$ws.Range("B7").Value = "MORON"
$ws.Range("C7").Value = "S1AGL065"

# Row 8: CP007_Promociones
$ws.Range("A8").Value = "CP007_Promociones"
$ws.Range("B8").Value = "5 ELEMENTOS BAR"

# Row 9: CP008_Promociones2
$ws.Range("A9").Value = "CP008_Promociones2"
$ws.Range("B9").Value = "A LA PIPETUA"

Write-Output "values set"
